# The deck ships two DrawingML theme parts: ppt/theme/theme1.xml ("Office
# Theme" colours) and ppt/theme/theme2.xml ("Integral" colours) - the slide
# master (and therefore every slide) is wired to the "Integral" theme.
#
# The authored change swaps the two themes' colour schemes (and names) so
# that the deck's live/visible theme becomes "Office Theme" colours while
# the other part keeps the "Integral" palette. The only theme surface the
# PowerPoint object model exposes for editing is the live/active theme
# (reached equivalently through Slide.ThemeColorScheme,
# SlideMaster.Theme.ThemeColorScheme, etc.) so we repaint its twelve theme
# colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) with the
# "Office Theme" palette that should end up "in front".

function ConvertTo-BgrInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Target palette: the original "Office Theme" (ppt/theme/theme1.xml) colours,
# now promoted to the live theme.
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-BgrInt $officeThemeColors[$i - 1]
}
